# Apply a price update pass over the cat-food price list:
#  - Rows 30-45 (column B) get new, higher prices.
#  - Row 40 additionally swaps its product name from the discontinued
#    "Sieger-Lata-Ad-Extra" to the new "Sieger-Lata-Ad-Recovery".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New price list for B30:B45 (in row order).
$newPrices = @{
    30 = 900
    31 = 900
    32 = 900
    33 = 900
    34 = 900
    35 = 3880
    36 = 3880
    37 = 3880
    38 = 2040
    39 = 5300
    40 = 5860
    41 = 2800
    42 = 2800
    43 = 1700
    44 = 1700
    45 = 1700
}

foreach ($row in $newPrices.Keys) {
    $ws.Cells.Item($row, 2).Value = $newPrices[$row]
}

# Product rename: "Sieger-Lata-Ad-Extra" -> "Sieger-Lata-Ad-Recovery"
$ws.Range("A40").Value = "Sieger-Lata-Ad-Recovery"
